# "borre el 2do hola de la oracion"
# The paragraph originally reads "Hola bebe como estas. Hola " (two runs:
# "Hola bebe como estas. " and "Hola "). The second "Hola " greeting is
# removed entirely, together with the space that used to separate it from
# the first sentence, leaving just "Hola bebe como estas."

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Hola bebe como estas. Hola ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Hola bebe como estas.", 2
) | Out-Null
